# Applies the "applied LOD to Tableau and finalized slides for DE" edit:
#  1) Slide 5 ("Tableau build"): append ", LOD" to the Tableau Functions bullet.
#  2) Slide 6 ("Tableau visualization"): turn "Link to Dashboard's" into
#     "Link to Dashboards: <hyperlink-url>" and add a trailing blank paragraph.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 5 - "Tableau Functions Used: Hover filter, Cluster Model (Analytics)"
#    -> "...Cluster Model (Analytics), LOD"
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$tr5 = $s5.Shapes.Item(2).TextFrame.TextRange

$oldBullet = "Tableau Functions Used: Hover filter, Cluster Model (Analytics)"
$newBullet = "Tableau Functions Used: Hover filter, Cluster Model (Analytics), LOD"

$full5 = $tr5.Text
$idx5 = $full5.IndexOf($oldBullet)
if ($idx5 -ge 0) {
    $range5 = $tr5.Characters($idx5 + 1, $oldBullet.Length)
    $range5.Text = $newBullet
}

# ---------------------------------------------------------------------------
# 2) Slide 6 - "Link to Dashboard's" -> "Link to Dashboards: <url>" (hyperlinked)
#    plus a new trailing empty paragraph.
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$tr6 = $s6.Shapes.Item(2).TextFrame.TextRange

$oldLink = "Link to Dashboard" + [char]8217 + "s"
$full6 = $tr6.Text
$idx6 = $full6.IndexOf($oldLink)
if ($idx6 -lt 0) {
    # fall back in case the apostrophe round-trips as a plain quote
    $oldLink = "Link to Dashboard's"
    $idx6 = $full6.IndexOf($oldLink)
}

if ($idx6 -ge 0) {
    # Rewrite the existing run's text (keeps it a single run, same as before).
    $range6 = $tr6.Characters($idx6 + 1, $oldLink.Length)
    $range6.Text = "Link to Dashboards"

    # Add the ": " separator as its own run.
    [void]$tr6.InsertAfter(": ")

    # Add the hyperlink display text as its own run.
    $url = "https://public.tableau.com/profile/dustin.elery#!/vizhome/SharkTankDashboard/DealsDashboard?publish=yes"
    [void]$tr6.InsertAfter($url)

    # Add a new trailing paragraph (inserted as "<CR><placeholder>", then the
    # placeholder character is deleted so the new paragraph is left empty).
    [void]$tr6.InsertAfter([char]13 + "#")
    $full6b = $tr6.Text
    $placeholder = $tr6.Characters($full6b.Length, 1)
    [void]$placeholder.Delete()

    # Turn the URL run into a real hyperlink.
    $full6c = $tr6.Text
    $urlIdx = $full6c.IndexOf($url)
    if ($urlIdx -ge 0) {
        $urlRange = $tr6.Characters($urlIdx + 1, $url.Length)
        $urlRange.ActionSettings(1).Hyperlink.Address = $url
    }
}
